$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.010.83'
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').Value = '2.487.51'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''323.01'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = '''104.48'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '''0.521'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.539'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('D10').Value = '''36.95'
$ws.Range('E10').Value = '  +2.61%  '
$ws.Range('D11').Value = '''0.0812'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '''18.25'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = '''7.19'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '2.870.74'
$ws.Range('E15').Value = '  +1.43%  '
$ws.Range('D16').Value = '2.448.57'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '''0.840'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '46.899.21'
$ws.Range('E18').Value = '  +3.89%  '
$ws.Range('D19').Value = '''12.54'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '''6.56'
$ws.Range('E20').Value = '  +2.72%  '
$ws.Range('D21').Value = '0.0₃0931'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').Value = '''70.74'
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').Value = '''250.07'
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('D24').Value = '''2.35'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').Value = '''26.08'
$ws.Range('E26').Value = '  +1.99%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').Value = '''10.06'
$ws.Range('E28').Value = '  +4.59%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '''34.87'
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('D31').Value = '''0.133'
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('D32').Value = '''49.52'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').Value = '''19.59'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('D34').Value = '''5.30'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '''0.0776'
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '''1.92'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = '''4.56'
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').Value = '''2.96'
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '''0.111'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''121.94'
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('D43').Value = '''21.53'
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').Value = '1.949.54'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '''2.97'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('D47').Value = '''2.11'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '''5.34'
$ws.Range('E50').Value = '  +12.88%  '
$ws.Range('D51').Value = '''78.76'
$ws.Range('E51').Value = '  +3.25%  '
